$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Locate the paragraph that ends with "...if she stopped at 1000. " (the
# last sentence of Problem 3's description) so we can anchor the new
# paragraphs right after it, regardless of its absolute index.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$anchorIndex = 0
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "counted to 10, and if she stopped at 100, and if she stopped at 1000\.") {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq 0) {
    throw "Could not find anchor sentence for Problem 3"
}
$A = $anchorIndex

# NOTE: paragraph-object handles become stale once earlier content in the
# document shifts around them, so every paragraph below is re-fetched by
# its (stable) numeric index right before it is used instead of being
# cached across mutating calls.

# ---------------------------------------------------------------------
# 1) Blank spacer paragraph, then the "So there is a pattern here..."
#    paragraph.
# ---------------------------------------------------------------------
$d.Paragraphs.Item($A).Range.InsertParagraphAfter()
$d.Paragraphs.Item($A + 1).Range.InsertAfter("So there is a pattern here. Now the idea I think would be to understand the pattern and then to be able to predict it. Also in the question it give the answer to A, " + [char]34 + [char]0x2026 + "after which she calls her first finger 10 and so on." + [char]34 + " ")

# ---------------------------------------------------------------------
# 2) Blank spacer paragraph, then the "I think the overall goal..."
#    paragraph, built from two separate runs.
# ---------------------------------------------------------------------
$d.Paragraphs.Item($A + 1).Range.InsertParagraphAfter()
$d.Paragraphs.Item($A + 2).Range.InsertParagraphAfter()
$d.Paragraphs.Item($A + 3).Range.InsertAfter("I think the overall go")
$d.Paragraphs.Item($A + 3).Range.InsertAfter("al is to be able to predict on what finger she will end up by following the pattern.")

# ---------------------------------------------------------------------
# 3) Blank spacer paragraph right before the pre-existing bold/bookmark
#    paragraph.
# ---------------------------------------------------------------------
$d.Paragraphs.Item($A + 3).Range.InsertParagraphAfter()

# ---------------------------------------------------------------------
# 4) The pre-existing empty paragraph (bold paragraph mark + _GoBack
#    bookmark) now sits at $A + 5. Split its bold paragraph-mark
#    formatting off into a brand-new trailing paragraph, then clear the
#    bold from the original (bookmarked) paragraph and give it a single
#    new run containing one space.
# ---------------------------------------------------------------------
$bmIndex = $A + 5

# 4a) Split: inserting a paragraph break after the still-bold paragraph
#     makes the new, following paragraph inherit the bold paragraph
#     mark (plus a throwaway empty run we clean up next).
$d.Paragraphs.Item($bmIndex).Range.InsertParagraphAfter()

$newBoldIndex = $bmIndex + 1
$d.Paragraphs.Item($newBoldIndex).Range.InsertAfter("X")
$tmpRange = $d.Paragraphs.Item($newBoldIndex).Range
$d.Range($tmpRange.Start, $tmpRange.Start + 1).Delete()

# 4b) Clear the bold paragraph mark on the original bookmark paragraph
#     (a zero-length range's formatting can't be toggled directly, so
#     round-trip a throwaway character through it).
$d.Paragraphs.Item($bmIndex).Range.InsertAfter("X")
$d.Paragraphs.Item($bmIndex).Range.Font.Bold = $false
$tmpRange2 = $d.Paragraphs.Item($bmIndex).Range
$d.Range($tmpRange2.Start, $tmpRange2.Start + 1).Delete()

# 4c) Finally, give the bookmark paragraph its new single-space run.
$d.Paragraphs.Item($bmIndex).Range.InsertAfter(" ")
